# Klas heeft lijst van leerlingen met voor- en achternaam:
# append each student's last name to the existing first-name cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Jan, Protput"
$ws.Range("C1").Value = "Piet, De Hans"
$ws.Range("D1").Value = "Nelis, Cornelis"
$ws.Range("E1").Value = "Corneel, Teeuwen"

$ws.Range("B2").Value = "Joris, Boris"
$ws.Range("C2").Value = "Welsey, De Kleine"
$ws.Range("D2").Value = "Diaby, Abdoulay"
$ws.Range("E2").Value = "Thomas, Vermeel"

# widen the (now longer) B column to fit its new content
$ws.Columns("B:B").AutoFit() | Out-Null

# last active cell ends up on E3
$ws.Range("E3").Select() | Out-Null
